$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.211.92'
$ws.Range("E2").Value = '  -4.60%  '
$ws.Range("D3").Value = '2.982.39'
$ws.Range("E3").Value = '  -5.93%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.99'
$ws.Range("E5").Value = '  -3.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '124.43'
$ws.Range("E6").Value = '  -7.72%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '2.979.88'
$ws.Range("E8").Value = '  -5.92%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.497'
$ws.Range("E9").Value = '  -3.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.133'
$ws.Range("E10").Value = '  -6.26%  '
$ws.Range("E11").Value = '  -3.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.434'
$ws.Range("E12").Value = '  -4.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000222'
$ws.Range("E13").Value = '  -6.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.43'
$ws.Range("E14").Value = '  -7.24%  '
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").Value = '3.478.24'
$ws.Range("E16").Value = '  -5.87%  '
$ws.Range("D17").Value = '60.350.32'
$ws.Range("E17").Value = '  -4.34%  '
$ws.Range("D18").Value = '2.989.13'
$ws.Range("E18").Value = '  -5.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.13'
$ws.Range("E19").Value = '  -7.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '429.15'
$ws.Range("E20").Value = '  -7.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.00'
$ws.Range("E21").Value = '  -6.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.657'
$ws.Range("E22").Value = '  -6.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.12'
$ws.Range("E23").Value = '  -6.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.86'
$ws.Range("E24").Value = '  -4.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.78'
$ws.Range("E25").Value = '  -5.42%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("E28").Value = '  -6.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.08'
$ws.Range("E29").Value = '  -8.49%  '
$ws.Range("E30").Value = '  -8.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.10'
$ws.Range("E31").Value = '  -7.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.97'
$ws.Range("E32").Value = '  -11.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0922'
$ws.Range("E33").Value = '  -10.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.23'
$ws.Range("E34").Value = '  -5.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.938'
$ws.Range("E35").Value = '  -9.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.51'
$ws.Range("E36").Value = '  -5.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '49.35'
$ws.Range("E37").Value = '  -3.99%  '
$ws.Range("D38").Value = '0.0₃0645'
$ws.Range("E38").Value = '  -8.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0356'
$ws.Range("E39").Value = '  -8.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.81'
$ws.Range("E40").Value = '  -3.95%  '
$ws.Range("E41").Value = '  -5.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '370.19'
$ws.Range("E42").Value = '  -8.68%  '
$ws.Range("D43").Value = '2.647.47'
$ws.Range("E43").Value = '  -5.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.36'
$ws.Range("E44").Value = '  -8.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.232'
$ws.Range("E46").Value = '  -7.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.88'
$ws.Range("E47").Value = '  -4.16%  '
$ws.Range("E48").Value = '  -8.86%  '
$ws.Range("E49").Value = '  -5.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.15'
$ws.Range("E50").Value = '  -8.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '31.23'
$ws.Range("E51").Value = '  -7.97%  '
